$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Delete rows 4 through 11 (the "1.1." .. "3.2." rows), working from the
# bottom up so indices of earlier rows stay valid as we delete.
for ($i = 11; $i -ge 4; $i--) {
    $tbl.Rows.Item($i).Delete()
}
